$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 394 (existing rows 394:422 shift down to 396:424)
$ws.Rows("394:395").Insert()

# New row 394 values
$ws.Range("A394").Value = 9
$ws.Range("B394").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C394").Value = "Metropolitana"
$ws.Range("D394").Value = 45021
$ws.Range("E394").Value = 13
$ws.Range("F394").Value = 100112017
$ws.Range("G394").Value = "Apio"
$ws.Range("H394").Value = "Americana (o)"
$ws.Range("I394").Value = "Primera"
$ws.Range("J394").Value = 90
$ws.Range("K394").Value = 8000
$ws.Range("L394").Value = 9000
$ws.Range("M394").Value = 8556
$ws.Range("N394").Value = "`$/docena de matas"
$ws.Range("O394").Value = "Región de Coquimbo"
$ws.Range("P394").Value = 1426
$ws.Range("Q394").Value = 6
$ws.Range("R394").Value = "Hortaliza"

# New row 395 values
$ws.Range("A395").Value = 9
$ws.Range("B395").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C395").Value = "Metropolitana"
$ws.Range("D395").Value = 45021
$ws.Range("E395").Value = 13
$ws.Range("F395").Value = 100112017
$ws.Range("G395").Value = "Apio"
$ws.Range("H395").Value = "Americana (o)"
$ws.Range("I395").Value = "Segunda"
$ws.Range("J395").Value = 55
$ws.Range("K395").Value = 7000
$ws.Range("L395").Value = 7000
$ws.Range("M395").Value = 7000
$ws.Range("N395").Value = "`$/docena de matas"
$ws.Range("O395").Value = "Región de Coquimbo"
$ws.Range("P395").Value = 1167
$ws.Range("Q395").Value = 6
$ws.Range("R395").Value = "Hortaliza"
